# Update "想去人数" (interested-count) values in column F across sheets,
# matching the generated-data refresh captured in the commit diff.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F6").Value  = 360
$ws1.Range("F9").Value  = 1439
$ws1.Range("F12").Value = 2994
$ws1.Range("F13").Value = 395
$ws1.Range("F14").Value = 1606
$ws1.Range("F15").Value = 1359
$ws1.Range("F18").Value = 1370
$ws1.Range("F19").Value = 263
$ws1.Range("F21").Value = 1116
$ws1.Range("F23").Value = 3461

# Sheet "演出" (sheet2)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F5").Value  = 25
$ws2.Range("F12").Value = 79

# Sheet "全部类型" (sheet4)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value  = 25
$ws4.Range("F16").Value = 360
$ws4.Range("F19").Value = 1439
$ws4.Range("F22").Value = 2994
$ws4.Range("F23").Value = 395
$ws4.Range("F24").Value = 1606
$ws4.Range("F25").Value = 1359
$ws4.Range("F28").Value = 1370
$ws4.Range("F29").Value = 263
$ws4.Range("F33").Value = 1116
$ws4.Range("F35").Value = 3461
$ws4.Range("F39").Value = 79
